$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 152
$ws.Range("I12").Value = 122.4
$ws.Range("K12").Value = 122.4
$ws.Range("M12").Value = 47.59999999999999

# Row 132
$ws.Range("H132").Value = 2465.84
$ws.Range("I132").Value = 1376.0435
$ws.Range("K132").Value = 4128.1305
$ws.Range("M132").Value = -1598.1305

# Row 137
$ws.Range("H137").Value = 1498.4117
$ws.Range("I137").Value = 1198.8182
$ws.Range("K137").Value = 3596.4546
$ws.Range("M137").Value = -1046.4546

$ws = $wb.Worksheets.Item("ARM")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# Row 55
$ws.Range("H55").Value = 32999.332
$ws.Range("J55").Value = 32999.332
$ws.Range("L55").Value = 32999.332
$ws.Range("N55").Value = -33629.332

# Row 61
$ws.Range("H61").Value = 8693.223
$ws.Range("I61").Value = 8693.223
$ws.Range("K61").Value = 8693.223
$ws.Range("M61").Value = -8481.223

# Row 102
$ws.Range("H102").Value = 2700
$ws.Range("I102").Value = 1344.5
$ws.Range("K102").Value = 1344.5
$ws.Range("M102").Value = 277.5

# Row 110
$ws.Range("H110").Value = 2441.652
$ws.Range("I110").Value = 1396.5
$ws.Range("K110").Value = 1396.5
$ws.Range("M110").Value = 648.5

# Row 136
$ws.Range("H136").Value = 8693.223
$ws.Range("I136").Value = 8693.223
$ws.Range("K136").Value = 26079.669
$ws.Range("M136").Value = -23529.669

$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Range("H16").Value = 5800
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5800
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5800
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -6140

# Row 99
$ws.Range("H99").Value = 1897.4286
$ws.Range("J99").Value = 999.5
$ws.Range("L99").Value = 999.5
$ws.Range("N99").Value = -3995.5

# Row 105
$ws.Range("H105").Value = 2132.65
$ws.Range("I105").Value = 1477.909
$ws.Range("J105").Value = 2932.889
$ws.Range("K105").Value = 1477.909
$ws.Range("L105").Value = 2932.889
$ws.Range("M105").Value = 269.0909999999999
$ws.Range("N105").Value = -6426.889

# Row 107
$ws.Range("H107").Value = 937
$ws.Range("I107").Value = 950.5
$ws.Range("K107").Value = 950.5
$ws.Range("M107").Value = 969.5

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# Row 99
$ws.Range("H99").Value = 2854.5334
$ws.Range("I99").Value = 3022.4
$ws.Range("J99").Value = 2518.8
$ws.Range("K99").Value = 3022.4
$ws.Range("L99").Value = 2518.8
$ws.Range("M99").Value = -1524.4
$ws.Range("N99").Value = -5514.8

# Row 126
$ws.Range("H126").Value = 2854.5334
$ws.Range("I126").Value = 3022.4
$ws.Range("J126").Value = 2518.8
$ws.Range("K126").Value = 9067.200000000001
$ws.Range("L126").Value = 7556.400000000001
$ws.Range("M126").Value = -6597.200000000001
$ws.Range("N126").Value = -12496.4

# Row 132
$ws.Range("I132").Value = 2333.3333
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6999.999899999999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -4469.999899999999
$ws.Range("N132").Value = -14060

# Row 134
$ws.Range("H134").Value = 3272.5
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 33
$ws.Range("H33").Value = 62562.438
$ws.Range("I33").Value = 69
$ws.Range("J33").Value = 250042.75
$ws.Range("K33").Value = 414
$ws.Range("L33").Value = 1500256.5
$ws.Range("M33").Value = -131
$ws.Range("N33").Value = -1500822.5

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1995
$ws.Range("I97").Value = 1287.25
$ws.Range("J97").Value = 2509.7273
$ws.Range("K97").Value = 1287.25
$ws.Range("L97").Value = 2509.7273
$ws.Range("M97").Value = -791.25
$ws.Range("N97").Value = -3501.7273

$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 3908.3333
$ws.Range("I4").Value = 3859
$ws.Range("K4").Value = 3859
$ws.Range("M4").Value = -3746

# Row 16
$ws.Range("H16").Value = 3000
$ws.Range("J16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("N16").Value = -3340

# Row 18
$ws.Range("H18").Value = 27997.5
$ws.Range("J18").Value = 7995
$ws.Range("L18").Value = 7995
$ws.Range("N18").Value = -8339

# Row 22
$ws.Range("H22").Value = 3184.5
$ws.Range("J22").Value = 2250
$ws.Range("L22").Value = 2250
$ws.Range("N22").Value = -2840

# Row 27
$ws.Range("H27").Value = 3184.5
$ws.Range("J27").Value = 2250
$ws.Range("L27").Value = 2250
$ws.Range("N27").Value = -2464

# Row 28
$ws.Range("H28").Value = 3908.3333
$ws.Range("I28").Value = 3859
$ws.Range("K28").Value = 3859
$ws.Range("M28").Value = -3627

# Row 37
$ws.Range("H37").Value = 3908.3333
$ws.Range("I37").Value = 3859
$ws.Range("K37").Value = 3859
$ws.Range("M37").Value = -3752

# Row 40
$ws.Range("H40").Value = 5187.25
$ws.Range("I40").Value = 4999.6665
$ws.Range("K40").Value = 4999.6665
$ws.Range("M40").Value = -4863.6665

# Row 47
$ws.Range("H47").Value = 19497
$ws.Range("I47").Value = 19497
$ws.Range("K47").Value = 19497
$ws.Range("M47").Value = -19007

# Row 52
$ws.Range("H52").Value = 19497
$ws.Range("I52").Value = 19497
$ws.Range("K52").Value = 19497
$ws.Range("M52").Value = -19264

# Row 100
$ws.Range("H100").Value = 5649.4287
$ws.Range("I100").Value = 5649.2
$ws.Range("K100").Value = 5649.2
$ws.Range("M100").Value = -5108.2

# Row 122
$ws.Range("H122").Value = 4472.579
$ws.Range("I122").Value = 3082.0833
$ws.Range("J122").Value = 6856.2856
$ws.Range("K122").Value = 9246.249899999999
$ws.Range("L122").Value = 20568.8568
$ws.Range("M122").Value = -6796.249899999999
$ws.Range("N122").Value = -25468.8568

$ws = $wb.Worksheets.Item("WVR")
# Row 12
$ws.Range("H12").Value = 2899
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

# Row 22
$ws.Range("H22").Value = 12500
$ws.Range("I22").Value = 12500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 12500
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -12207
$ws.Range("N22").ClearContents()

# Row 40
$ws.Range("H40").Value = 32129.5
$ws.Range("J40").Value = 32129.5
$ws.Range("L40").Value = 32129.5
$ws.Range("N40").Value = -32427.5

# Row 100
$ws.Range("H100").Value = 1938.0454
$ws.Range("I100").Value = 2163.4285
$ws.Range("K100").Value = 4326.857
$ws.Range("M100").Value = -3785.857
